$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C8").Value = "RCK99"
$ws.Range("C9:C15").Value = "RCK100"
$ws.Range("C16:C22").Value = "RCK101"
$ws.Range("C23:C29").Value = "RCK102"
$ws.Range("C30:C36").Value = "RCK103"
$ws.Range("C37:C43").Value = "RCK104"
